$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.884.00'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '2.215.54'
$ws.Range('E3').Value = '  -1.63%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '257.49'
$ws.Range('E5').Value = '  +5.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.618'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '76.75'
$ws.Range('E7').Value = '  +1.34%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.595'
$ws.Range('E9').Value = '  -1.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.02'
$ws.Range('E10').Value = '  +2.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0909'
$ws.Range('E11').Value = '  -2.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.96'
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('D14').Value = '2.545.54'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.50'
$ws.Range('E15').Value = '  -1.39%  '
$ws.Range('D16').Value = '2.209.98'
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.784'
$ws.Range('E17').Value = '  -2.06%  '
$ws.Range('D18').Value = '42.818.03'
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.17'
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.69'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.20'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.32'
$ws.Range('E24').Value = '  -6.98%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '43.25'
$ws.Range('E26').Value = '  +10.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.76'
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.35'
$ws.Range('E28').Value = '  -2.88%  '
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.18'
$ws.Range('E30').Value = '  -1.68%  '
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.37'
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0874'
$ws.Range('E33').Value = '  +9.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.23'
$ws.Range('E34').Value = '  -1.87%  '
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0362'
$ws.Range('E36').Value = '  +7.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.109'
$ws.Range('E37').Value = '  -1.42%  '
$ws.Range('E38').Value = '  -0.56%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.78'
$ws.Range('E39').Value = '  -2.99%  '
$ws.Range('E40').Value = '  -0.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.81'
$ws.Range('E41').Value = '  +17.26%  '
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.31'
$ws.Range('E43').Value = '  -3.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '60.06'
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.83'
$ws.Range('E45').Value = '  -2.84%  '
$ws.Range('E46').Value = '  -4.45%  '
$ws.Range('E47').Value = '  -1.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.464'
$ws.Range('E48').Value = '  -5.17%  '
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('E50').Value = '  -1.29%  '
$ws.Range('D51').Value = '2.433.55'
$ws.Range('E51').Value = '  -1.07%  '
